$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 8.91390722595512
$ws.Range("E2").Value = 8.423642158508301
$ws.Range("F2").Value = 9.261549913802089
$ws.Range("G2").Value = 7.4431128737111
$ws.Range("H2").Value = 86406484
$ws.Range("I2").Value = "YRD"

$ws.Range("D3").Value = 8.272106026459451
$ws.Range("E3").Value = 6.079285144805908
$ws.Range("F3").Value = 8.459298120327766
$ws.Range("G3").Value = 5.651417562399507
$ws.Range("H3").Value = 86406484
$ws.Range("I3").Value = "YRD"

$ws.Range("D4").Value = 10.09054264223088
$ws.Range("E4").Value = 10.9997615814209
$ws.Range("F4").Value = 13.25497985691706
$ws.Range("G4").Value = 8.985218357412636
$ws.Range("H4").Value = 86406484
$ws.Range("I4").Value = "YRD"

$ws.Range("D5").Value = 12.23879391921309
$ws.Range("E5").Value = 22.85525703430176
$ws.Range("F5").Value = 23.17615823078978
$ws.Range("G5").Value = 12.23879391921309
$ws.Range("H5").Value = 86406484
$ws.Range("I5").Value = "YRD"

$ws.Range("D6").Value = 21.17944484197945
$ws.Range("E6").Value = 25.85033226013184
$ws.Range("F6").Value = 28.52450456290409
$ws.Range("G6").Value = 20.88528595668222
$ws.Range("H6").Value = 86406484
$ws.Range("I6").Value = "YRD"

$ws.Range("D7").Value = 19.30752260255097
$ws.Range("E7").Value = 18.83508491516113
$ws.Range("F7").Value = 20.67135099809254
$ws.Range("G7").Value = 17.69410589685035
$ws.Range("H7").Value = 86406484
$ws.Range("I7").Value = "YRD"

$ws.Range("D8").Value = 22.46304836654575
$ws.Range("E8").Value = 22.24911308288575
$ws.Range("F8").Value = 23.31878270041382
$ws.Range("G8").Value = 20.66243858045711
$ws.Range("H8").Value = 86406484
$ws.Range("I8").Value = "YRD"

$ws.Range("D9").Value = 22.24911101845287
$ws.Range("E9").Value = 34.10460662841797
$ws.Range("F9").Value = 34.30962948141222
$ws.Range("G9").Value = 22.07974800987682
$ws.Range("H9").Value = 86406484
$ws.Range("I9").Value = "YRD"

$ws.Range("D10").Value = 39.2914544267158
$ws.Range("E10").Value = 40.02009963989258
$ws.Range("F10").Value = 49.34490345360663
$ws.Range("G10").Value = 37.1700854492182
$ws.Range("H10").Value = 86406484
$ws.Range("I10").Value = "YRD"

$ws.Range("D11").Value = 40.39826182526441
$ws.Range("E11").Value = 37.28076934814453
$ws.Range("F11").Value = 44.20750537308548
$ws.Range("G11").Value = 34.46764903097857
$ws.Range("H11").Value = 86406484
$ws.Range("I11").Value = "YRD"

$ws.Range("D12").Value = 36.89338499920883
$ws.Range("E12").Value = 32.53074264526367
$ws.Range("F12").Value = 37.4652314815369
$ws.Range("G12").Value = 31.82976635995937
$ws.Range("H12").Value = 86406484
$ws.Range("I12").Value = "YRD"

$ws.Range("D13").Value = 19.13965088347417
$ws.Range("E13").Value = 16.68561172485352
$ws.Range("F13").Value = 20.91511130096047
$ws.Range("G13").Value = 14.90085548958505
$ws.Range("H13").Value = 86406484
$ws.Range("I13").Value = "YRD"

$ws.Range("D14").Value = 17.34559972940203
$ws.Range("E14").Value = 14.46396255493164
$ws.Range("F14").Value = 17.71277657234554
$ws.Range("G14").Value = 12.78145756857307
$ws.Range("H14").Value = 86406484
$ws.Range("I14").Value = "YRD"

$ws.Range("D15").Value = 9.69531742521097
$ws.Range("E15").Value = 10.38319206237793
$ws.Range("F15").Value = 11.39641326694719
$ws.Range("G15").Value = 9.504757571489016
$ws.Range("H15").Value = 86406484
$ws.Range("I15").Value = "YRD"

$ws.Range("D16").Value = 12.93948301086654
$ws.Range("E16").Value = 14.72423934936523
$ws.Range("F16").Value = 15.44929722134772
$ws.Range("G16").Value = 12.10287886503816
$ws.Range("H16").Value = 86406484
$ws.Range("I16").Value = "YRD"

$ws.Range("D17").Value = 13.22764603742721
$ws.Range("E17").Value = 10.40178203582764
$ws.Range("F17").Value = 13.67290562047559
$ws.Range("G17").Value = 9.974184152654766
$ws.Range("H17").Value = 86406484
$ws.Range("I17").Value = "YRD"

$ws.Range("D18").Value = 6.525514949835163
$ws.Range("E18").Value = 6.014256477355957
$ws.Range("F18").Value = 6.971704025176414
$ws.Range("G18").Value = 5.614545606148818
$ws.Range("H18").Value = 86406484
$ws.Range("I18").Value = "YRD"

$ws.Range("D19").Value = 5.363564246385358
$ws.Range("E19").Value = 4.257387161254883
$ws.Range("F19").Value = 6.860156825473529
$ws.Range("G19").Value = 4.229500593543048
$ws.Range("H19").Value = 86406484
$ws.Range("I19").Value = "YRD"

$ws.Range("D20").Value = 3.718242358941158
$ws.Range("E20").Value = 3.671764373779297
$ws.Range("F20").Value = 4.043588476698703
$ws.Range("G20").Value = 3.05825430476918
$ws.Range("H20").Value = 86406484
$ws.Range("I20").Value = "YRD"

$ws.Range("D21").Value = 3.774015674865765
$ws.Range("E21").Value = 3.337122201919556
$ws.Range("F21").Value = 4.164431165929292
$ws.Range("G21").Value = 3.253461834530728
$ws.Range("H21").Value = 86406484
$ws.Range("I21").Value = "YRD"

$ws.Range("D22").Value = 2.74220351393124
$ws.Range("E22").Value = 3.05825400352478
$ws.Range("F22").Value = 3.950632117229474
$ws.Range("G22").Value = 2.710598309834741
$ws.Range("H22").Value = 86406484
$ws.Range("I22").Value = "YRD"

$ws.Range("D23").Value = 3.151210330881089
$ws.Range("E23").Value = 3.21627950668335
$ws.Range("F23").Value = 3.606695004745919
$ws.Range("G23").Value = 2.82586400862078
$ws.Range("H23").Value = 86406484
$ws.Range("I23").Value = "YRD"

$ws.Range("D24").Value = 4.387525650473951
$ws.Range("E24").Value = 3.578808069229126
$ws.Range("F24").Value = 4.60132459711616
$ws.Range("G24").Value = 3.290644348667718
$ws.Range("H24").Value = 86406484
$ws.Range("I24").Value = "YRD"

$ws.Range("D25").Value = 5.484407166886446
$ws.Range("E25").Value = 3.774015665054321
$ws.Range("F25").Value = 5.586658948387019
$ws.Range("G25").Value = 3.662468508825878
$ws.Range("H25").Value = 86406484
$ws.Range("I25").Value = "YRD"

$ws.Range("D26").Value = 2.751499186378496
$ws.Range("E26").Value = 2.918820142745972
$ws.Range("F26").Value = 4.155135347067074
$ws.Range("G26").Value = 2.695725608130837
$ws.Range("H26").Value = 86406484
$ws.Range("I26").Value = "YRD"

$ws.Range("D27").Value = 2.649247589288955
$ws.Range("E27").Value = 2.286719083786011
$ws.Range("F27").Value = 2.9281159432137
$ws.Range("G27").Value = 2.093370532742806
$ws.Range("H27").Value = 86406484
$ws.Range("I27").Value = "YRD"

$ws.Range("D28").Value = 2.416857384603992
$ws.Range("E28").Value = 1.905599117279053
$ws.Range("F28").Value = 2.529334261037852
$ws.Range("G28").Value = 1.766165054446093
$ws.Range("H28").Value = 86406484
$ws.Range("I28").Value = "YRD"

$ws.Range("D29").Value = 1.636026510674624
$ws.Range("E29").Value = 1.199133038520813
$ws.Range("F29").Value = 1.736419084333834
$ws.Range("G29").Value = 1.161950653082105
$ws.Range("H29").Value = 86406484
$ws.Range("I29").Value = "YRD"

$ws.Range("D30").Value = 0.9667429916273974
$ws.Range("E30").Value = 0.7501553893089294
$ws.Range("F30").Value = 1.003925491098477
$ws.Range("G30").Value = 0.725057243720548
$ws.Range("H30").Value = 86406484
$ws.Range("I30").Value = "YRD"

$ws.Range("D31").Value = 1.264202290810901
$ws.Range("E31").Value = 3.05825400352478
$ws.Range("F31").Value = 3.17909697466685
$ws.Range("G31").Value = 1.236315502460691
$ws.Range("H31").Value = 86406484
$ws.Range("I31").Value = "YRD"

$ws.Range("D32").Value = 1.952077122086589
$ws.Range("E32").Value = 1.970668315887451
$ws.Range("F32").Value = 2.258832263049833
$ws.Range("G32").Value = 1.812643168580125
$ws.Range("H32").Value = 86406484
$ws.Range("I32").Value = "YRD"

$ws.Range("D33").Value = 2.15658045336517
$ws.Range("E33").Value = 2.491222381591797
$ws.Range("F33").Value = 2.584178349641132
$ws.Range("G33").Value = 1.979964114071432
$ws.Range("H33").Value = 86406484
$ws.Range("I33").Value = "YRD"

$ws.Range("D34").Value = 2.277423466652015
$ws.Range("E34").Value = 2.054328918457031
$ws.Range("F34").Value = 2.31460585441376
$ws.Range("G34").Value = 2.00785093375485
$ws.Range("H34").Value = 86406484
$ws.Range("I34").Value = "YRD"

$ws.Range("D35").Value = 2.900228855288518
$ws.Range("E35").Value = 3.997110605239868
$ws.Range("F35").Value = 4.64780282318635
$ws.Range("G35").Value = 2.807272887474594
$ws.Range("H35").Value = 86406484
$ws.Range("I35").Value = "YRD"

$ws.Range("D36").Value = 4.294569656883831
$ws.Range("E36").Value = 4.684985160827637
$ws.Range("F36").Value = 5.084696039830145
$ws.Range("G36").Value = 4.229500701725873
$ws.Range("H36").Value = 86406484
$ws.Range("I36").Value = "YRD"

$ws.Range("D37").Value = 4.275978260734711
$ws.Range("E37").Value = 4.415412425994873
$ws.Range("F37").Value = 4.63850673581194
$ws.Range("G37").Value = 4.062179325535525
$ws.Range("H37").Value = 86406484
$ws.Range("I37").Value = "YRD"

$ws.Range("D38").Value = 5.626877779533838
$ws.Range("E38").Value = 5.452850341796875
$ws.Range("F38").Value = 8.894720441853881
$ws.Range("G38").Value = 5.143469049845812
$ws.Range("H38").Value = 86406484
$ws.Range("I38").Value = "YRD"

$ws.Range("D39").Value = 4.544042091910759
$ws.Range("E39").Value = 5.491523265838623
$ws.Range("F39").Value = 5.684887043567184
$ws.Range("G39").Value = 4.418356166553883
$ws.Range("H39").Value = 86406484
$ws.Range("I39").Value = "YRD"

$ws.Range("D40").Value = 6.700045310009397
$ws.Range("E40").Value = 6.322986602783203
$ws.Range("F40").Value = 7.193122471086708
$ws.Range("G40").Value = 4.544042201897327
$ws.Range("H40").Value = 86406484
$ws.Range("I40").Value = "YRD"

$ws.Range("D41").Value = 5.909999847412109
$ws.Range("E41").Value = 5.920000076293945
$ws.Range("F41").Value = 6.701000213623047
$ws.Range("G41").Value = 5.849999904632568
$ws.Range("H41").Value = 86406484
$ws.Range("I41").Value = "YRD"

